# bayiler bulk tamamlandi. file düzenlendi
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clean up / normalize the "Bölge Şehir" text blocks: drop stray trailing
# separators left over from the bulk dealer edit, and fix the spacing that
# crept into the Zafer GENÇ contact entry. (Order matches the author's
# edit sequence.)
$ws.Range("I5").Value = "İSTANBUL=ÜMRANİYE,ÜSKÜDAR,ATAŞEHİR,BEYKOZ,KADIKÖY"
$ws.Range("I10").Value = "GAZİANTEP=ARABAN,KARKAMIŞ,NİZİP,OĞUZELİ,ŞAHİNBEY,ŞEHİTKAMİL,YAVUZELİ"
$ws.Range("I11").Value = "HATAY=ARSUZ,BELEN,DÖRTYOL,İSKENDERUN"
$ws.Range("I13").Value = "HATAY=ALTINÖZÜ,ANTAKYA,BELEN,DEFNE,HASSA,KIRIKHAN,KUMLU,MERKEZ,REYHANLI,SAMANDAĞ,YAYLADAĞI"
$ws.Range("I14").Value = "ŞANLIURFA=AKÇAKALE,BİRECİK,BOZOVA,CEYLANPINAR,EYYÜBİYE,HALFETİ,HALİLİYE,HARRAN,HİLVAN,KARAKÖPRÜ,MERKEZ,SİVEREK,VİRANŞEHİR"
$ws.Range("I2").Value = "İSTANBUL=KARTAL,KADIKÖY,MALTEPE,PENDİK,TUZLA,ADALAR;`nKOCAELİ=DARICA,GEBZE;"
$ws.Range("H2").Value = " Zafer GENÇ < zafergenc02@gmail.com> "

# Restore the active selection that was left on the sheet (H2 selected,
# with the window scrolled right so column E is the first visible column).
$ws.Range("H2").Select()
